# Update "想去人数" (want-to-go count) figures in column F for both the
# "展览" (Exhibition) and "全部类型" (All Types) sheets, which hold
# identical data tables.

$wb = $excel.ActiveWorkbook

# Row -> new F-column value (per the refreshed scrape).
$updates = @{
    3  = 7132   # 合肥·第二届漫画城市动漫展 -故事再次开始   (7112 -> 7132)
    4  = 5126   # 合肥·环形宇宙动漫游戏嘉年华                (5094 -> 5126)
    6  = 162    # ...内场票-钱文青                            (159  -> 162)
    9  = 102    # ...内场-《球声雅集》                        (96   -> 102)
    11 = 82     # 合肥·首届运动番only                         (79   -> 82)
    13 = 628    # 合肥·Look Look动漫嘉年华                    (623  -> 628)
    14 = 184    # 合肥·第十三届次元之门动漫游戏博览会          (171  -> 184)
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
